$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new log row 13: 11/4/2023, "~3 hrs", translator scanning work ---
$ws.Range("A13").Value = (Get-Date -Year 2023 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B13").Value = "~3 hrs"
$ws.Range("C13").Value = "Completed the translator scanning portion and tested it using a temp printing function. Corrected some bugs related to the scanning."
$ws.Range("C13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 60

# --- Add new row 14: just the date, rest left blank for future entries ---
$ws.Range("A14").Value = (Get-Date -Year 2023 -Month 11 -Day 5 -Hour 0 -Minute 0 -Second 0)

# --- View: freeze the header row and scroll/select near the new rows ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B14").Select()
